# Adds a new "Misc" row (MDM.PATH_ABS) to the options table, and updates
# the sheet view/selection to reflect where the author ended up after
# finishing the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new row 29 data -----------------------------------------
# Set cells in A, D, E, B, C order so that new shared strings are appended
# in the same order the original author's save produced them:
#   77 "Misc", 78 "Old BV requires ...", 79 "MDM.PATH_ABS",
#   80 "Use absolute paths in multi-participant MDMs"
$ws.Range("A29").Value = "Misc"
$ws.Range("D29").Value = "Old BV requires absolute paths, but relative paths are more flexible (can move the directory, use on another computer, etc)."
$ws.Range("E29").Value = "MDM.PATH_ABS"
$ws.Range("B29").Value = "Use absolute paths in multi-participant MDMs"
$ws.Range("C29").Value = $false

# Row 29 needs to be taller to accommodate the wrapped description text.
$ws.Rows.Item(29).RowHeight = 30

# --- Grow the table / autofilter to include the new row ------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E29"))

# --- Update the view: scrolled up a bit, with B29 now selected -----------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B29").Select() | Out-Null
